$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 14-15: continuation of the Priority/Difficulty table (cols C/D/E)
$ws.Range("C14").Value = "Changing Tools"
$ws.Range("D14").Value = "Medium "
$ws.Range("E14").Value = "Easy"

$ws.Range("C15").Value = "Upload data to graph"
$ws.Range("D15").Value = "High "
$ws.Range("E15").Value = "Medium"

# New table starting row 18: Person / Tasks
$ws.Range("A18").Value = "Person"
$ws.Range("B18").Value = "Tasks"

$ws.Range("A19").Value = "Kat "
$ws.Range("B19").Value = "Delete Lines "

$ws.Range("B20").Value = "Colored Boxes"

$ws.Range("B21").Value = "Uploading data "

$ws.Range("A22").Value = "Fan "
$ws.Range("B22").Value = "Changing tools "

$ws.Range("B23").Value = "add segments to existing lines"

$ws.Range("A24").Value = "Alex"
$ws.Range("B24").Value = "Lightweight STL editor"

# Column B is slightly wider now (one extra character of width)
$ws.Columns.Item(2).ColumnWidth = 23.83

# Selection becomes "select all" (whole sheet)
$ws.Cells.Select() | Out-Null
